# Insert a new data row before row 179 (shifts existing rows 179..215 down to 180..216)
# and populate the newly inserted row 179 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(179).Insert()

$ws.Cells.Item(179,1).Value  = 9
$ws.Cells.Item(179,2).Value  = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(179,3).Value  = 'Metropolitana'
$ws.Cells.Item(179,4).Value  = 44711
$ws.Cells.Item(179,5).Value  = 13
$ws.Cells.Item(179,6).Value  = 100112026
$ws.Cells.Item(179,7).Value  = 'Haba'
$ws.Cells.Item(179,8).Value  = 'Sin especificar'
$ws.Cells.Item(179,9).Value  = 'Primera'
$ws.Cells.Item(179,10).Value = 52
$ws.Cells.Item(179,11).Value = 20000
$ws.Cells.Item(179,12).Value = 21000
$ws.Cells.Item(179,13).Value = 20500
$ws.Cells.Item(179,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(179,15).Value = 'Región Metropolitana'
$ws.Cells.Item(179,16).Value = 820
$ws.Cells.Item(179,17).Value = 25
$ws.Cells.Item(179,18).Value = 'Hortaliza'
